$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: "label" styled like the other header cells (copy format from E1)
$ws.Range("F1").Value = "label"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data cells F2:F7: "real"
$ws.Range("F2").Value = "real"
$ws.Range("F3").Value = "real"
$ws.Range("F4").Value = "real"
$ws.Range("F5").Value = "real"
$ws.Range("F6").Value = "real"
$ws.Range("F7").Value = "real"
